$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1177.1111
$ws.Range("I28").Value = 979.1579
$ws.Range("K28").Value = 979.1579
$ws.Range("M28").Value = -494.1579
$ws.Range("H40").Value = 4755.7417
$ws.Range("I40").Value = 6079.5415
$ws.Range("J40").Value = 3919.658
$ws.Range("K40").Value = 6079.5415
$ws.Range("L40").Value = 3919.658
$ws.Range("M40").Value = -5904.5415
$ws.Range("N40").Value = -4269.657999999999
$ws.Range("H96").Value = 554.1818
$ws.Range("I96").Value = 554.1818
$ws.Range("K96").Value = 1662.5454
$ws.Range("M96").Value = -289.5454
$ws.Range("H106").Value = 11858.417
$ws.Range("I106").Value = 14206.353
$ws.Range("J106").Value = 6156.2856
$ws.Range("K106").Value = 14206.353
$ws.Range("L106").Value = 6156.2856
$ws.Range("M106").Value = -13575.353
$ws.Range("N106").Value = -7418.2856
$ws.Range("H131").Value = 6278.5454
$ws.Range("I131").Value = 3897.6
$ws.Range("K131").Value = 11692.8
$ws.Range("M131").Value = -6652.799999999999
$ws.Range("H132").Value = 3100.5522
$ws.Range("I132").Value = 3104.15
$ws.Range("J132").Value = 3069.7144
$ws.Range("K132").Value = 9312.450000000001
$ws.Range("L132").Value = 9209.143199999999
$ws.Range("M132").Value = -6782.450000000001
$ws.Range("N132").Value = -14269.1432
$ws.Range("H137").Value = 223134.38
$ws.Range("I137").Value = 223134.38
$ws.Range("K137").Value = 669403.14
$ws.Range("M137").Value = -666853.14
$ws.Range("H139").Value = 91999.8
$ws.Range("J139").Value = 91999.8
$ws.Range("L139").Value = 91999.8
$ws.Range("N139").Value = -102279.8
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 291.66666
$ws.Range("I4").Value = 292.2
$ws.Range("J4").Value = 289
$ws.Range("K4").Value = 292.2
$ws.Range("L4").Value = 289
$ws.Range("M4").Value = -176.2
$ws.Range("N4").Value = -521
$ws.Range("H11").Value = 5299.8
$ws.Range("J11").Value = 4833.3335
$ws.Range("L11").Value = 4833.3335
$ws.Range("N11").Value = -5121.3335
$ws.Range("H13").Value = 2316.3333
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 2316.3333
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 2316.3333
$ws.Range("M13").ClearContents()
$ws.Range("N13").Value = -2604.3333
$ws.Range("H61").Value = 3588.2727
$ws.Range("I61").Value = 3588.2727
$ws.Range("K61").Value = 3588.2727
$ws.Range("M61").Value = -3376.2727
$ws.Range("H132").Value = 5775.5454
$ws.Range("I132").Value = 3916
$ws.Range("K132").Value = 11748
$ws.Range("M132").Value = -9218
$ws.Range("H136").Value = 3588.2727
$ws.Range("I136").Value = 3588.2727
$ws.Range("K136").Value = 10764.8181
$ws.Range("M136").Value = -8214.8181
$ws.Range("H139").Value = 180555.42
$ws.Range("J139").Value = 180555.42
$ws.Range("L139").Value = 180555.42
$ws.Range("N139").Value = -190835.42
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3492.348
$ws.Range("I20").Value = 3159.2354
$ws.Range("J20").Value = 4436.1665
$ws.Range("K20").Value = 3159.2354
$ws.Range("L20").Value = 4436.1665
$ws.Range("M20").Value = -2912.2354
$ws.Range("N20").Value = -4930.1665
$ws.Range("H80").Value = 331.5
$ws.Range("I80").Value = 453.33334
$ws.Range("J80").Value = 310
$ws.Range("K80").Value = 453.33334
$ws.Range("L80").Value = 310
$ws.Range("M80").Value = 544.66666
$ws.Range("N80").Value = -2306
$ws.Range("H83").Value = 331.5
$ws.Range("I83").Value = 453.33334
$ws.Range("J83").Value = 310
$ws.Range("K83").Value = 2266.6667
$ws.Range("L83").Value = 1550
$ws.Range("M83").Value = 2725.3333
$ws.Range("N83").Value = -11534
$ws.Range("H86").Value = 14210.389
$ws.Range("I86").Value = 20809.637
$ws.Range("J86").Value = 3840.1428
$ws.Range("K86").Value = 20809.637
$ws.Range("L86").Value = 3840.1428
$ws.Range("M86").Value = -19686.637
$ws.Range("N86").Value = -6086.1428
$ws.Range("H89").Value = 14210.389
$ws.Range("I89").Value = 20809.637
$ws.Range("J89").Value = 3840.1428
$ws.Range("K89").Value = 104048.185
$ws.Range("L89").Value = 19200.714
$ws.Range("M89").Value = -98432.185
$ws.Range("N89").Value = -30432.714
$ws.Range("H105").Value = 52019.9
$ws.Range("I105").Value = 73319.07000000001
$ws.Range("K105").Value = 73319.07000000001
$ws.Range("M105").Value = -71572.07000000001
$ws.Range("H107").Value = 2503.8572
$ws.Range("I107").Value = 2421.1667
$ws.Range("J107").Value = 3000
$ws.Range("K107").Value = 2421.1667
$ws.Range("L107").Value = 3000
$ws.Range("M107").Value = -501.1667000000002
$ws.Range("N107").Value = -6840
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 40829.85
$ws.Range("I31").Value = 2156.7144
$ws.Range("K31").Value = 2156.7144
$ws.Range("M31").Value = -1861.7144
$ws.Range("H34").Value = 40829.85
$ws.Range("I34").Value = 2156.7144
$ws.Range("K34").Value = 2156.7144
$ws.Range("M34").Value = -1954.7144
$ws.Range("H107").Value = 370.33334
$ws.Range("I107").Value = 144.4
$ws.Range("K107").Value = 144.4
$ws.Range("M107").Value = 1775.6
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 940.2353000000001
$ws.Range("I14").Value = 940.2353000000001
$ws.Range("K14").Value = 2820.7059
$ws.Range("M14").Value = -2647.7059
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("N21").ClearContents()
$ws.Range("H56").Value = 6644.3887
$ws.Range("I56").Value = 6644.3887
$ws.Range("K56").Value = 6644.3887
$ws.Range("M56").Value = -6114.3887
$ws.Range("H68").Value = 3041.2126
$ws.Range("J68").Value = 3218.3635
$ws.Range("L68").Value = 9655.0905
$ws.Range("N68").Value = -11277.0905
$ws.Range("H71").Value = 3041.2126
$ws.Range("J71").Value = 3218.3635
$ws.Range("L71").Value = 28965.2715
$ws.Range("N71").Value = -37077.2715
$ws.Range("H81").Value = 4782.577
$ws.Range("I81").Value = 721.53845
$ws.Range("K81").Value = 2164.61535
$ws.Range("M81").Value = -1041.61535
$ws.Range("H84").Value = 4782.577
$ws.Range("I84").Value = 721.53845
$ws.Range("K84").Value = 6493.84605
$ws.Range("M84").Value = -877.8460500000001
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 10500
$ws.Range("I44").Value = 10500
$ws.Range("K44").Value = 10500
$ws.Range("M44").Value = -9904
$ws.Range("H122").Value = 601465.6
$ws.Range("I122").Value = 749349.0600000001
$ws.Range("J122").Value = 9931.666999999999
$ws.Range("K122").Value = 2248047.18
$ws.Range("L122").Value = 29795.001
$ws.Range("M122").Value = -2245597.18
$ws.Range("N122").Value = -34695.001
$ws.Range("H132").Value = 5997.875
$ws.Range("I132").Value = 5998.8335
$ws.Range("J132").Value = 5995
$ws.Range("K132").Value = 17996.5005
$ws.Range("L132").Value = 17985
$ws.Range("M132").Value = -15466.5005
$ws.Range("N132").Value = -23045
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 601.2
$ws.Range("I16").Value = 601.2
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 601.2
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -431.2
$ws.Range("N16").ClearContents()
$ws.Range("H46").Value = 11430.737
$ws.Range("J46").Value = 11763.143
$ws.Range("L46").Value = 11763.143
$ws.Range("N46").Value = -12139.143
$ws.Range("H55").Value = 1252.2609
$ws.Range("I55").Value = 714
$ws.Range("K55").Value = 714
$ws.Range("M55").Value = -541
$ws.Range("H82").Value = 2981.3333
$ws.Range("I82").Value = 1462.6666
$ws.Range("J82").Value = 4500
$ws.Range("K82").Value = 1462.6666
$ws.Range("L82").Value = 4500
$ws.Range("M82").Value = -1101.6666
$ws.Range("N82").Value = -5222
$ws.Range("H85").Value = 2981.3333
$ws.Range("I85").Value = 1462.6666
$ws.Range("J85").Value = 4500
$ws.Range("K85").Value = 1462.6666
$ws.Range("L85").Value = 4500
$ws.Range("M85").Value = -214.6666
$ws.Range("N85").Value = -6996
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 10747.16
$ws.Range("I62").Value = 35332.668
$ws.Range("J62").Value = 7394.591
$ws.Range("K62").Value = 35332.668
$ws.Range("L62").Value = 7394.591
$ws.Range("M62").Value = -34708.668
$ws.Range("N62").Value = -8642.591
$ws.Range("H65").Value = 10747.16
$ws.Range("I65").Value = 35332.668
$ws.Range("J65").Value = 7394.591
$ws.Range("K65").Value = 176663.34
$ws.Range("L65").Value = 36972.955
$ws.Range("M65").Value = -173543.34
$ws.Range("N65").Value = -43212.955
$ws.Range("H81").Value = 1675
$ws.Range("I81").Value = 1628.5714
$ws.Range("K81").Value = 3257.1428
$ws.Range("M81").Value = -2196.1428
$ws.Range("H84").Value = 1675
$ws.Range("I84").Value = 1628.5714
$ws.Range("K84").Value = 16285.714
$ws.Range("M84").Value = -10981.714
$ws.Range("H132").Value = 122910.54
$ws.Range("I132").Value = 1777.7903
$ws.Range("J132").Value = 381884
$ws.Range("K132").Value = 5333.3709
$ws.Range("L132").Value = 1145652
$ws.Range("M132").Value = -2803.3709
$ws.Range("N132").Value = -1150712
